$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.9635149859347223
$ws.Range("D2").Value = 0.03284191414364557
$ws.Range("E2").Value = 0.3616638282185383
$ws.Range("F2").Value = 0.8103571534480807
$ws.Range("G2").Value = 0.002438768298994252
$ws.Range("K2").Value = 0.3794774899967308
$ws.Range("L2").Value = 0.1096003322991379
$ws.Range("M2").Value = 0.1927093455183346
$ws.Range("N2").Value = 2.147324177306682
$ws.Range("O2").Value = 2.821822959616128

$ws.Range("B3").Value = 0.9448167195582187
$ws.Range("D3").Value = 0.03105887091933113
$ws.Range("E3").Value = 0.364042379067913
$ws.Range("F3").Value = 0.8069291697227783
$ws.Range("G3").Value = 0.002441190921033108
$ws.Range("K3").Value = 0.3326898033872112
$ws.Range("L3").Value = 0.1038290825171799
$ws.Range("M3").Value = 0.1873183987160267
$ws.Range("N3").Value = 2.167958780013794
$ws.Range("O3").Value = 2.823792047177506

$ws.Range("B4").Value = 0.9337688813018588
$ws.Range("D4").Value = 0.02995268740206569
$ws.Range("E4").Value = 0.3655912990694961
$ws.Range("F4").Value = 0.8052877644600258
$ws.Range("G4").Value = 0.002442759368287969
$ws.Range("K4").Value = 0.3038275644057364
$ws.Range("L4").Value = 0.1003331979103805
$ws.Range("M4").Value = 0.1840965728485386
$ws.Range("N4").Value = 2.181279997314441
$ws.Range("O4").Value = 2.826517546788523

$ws.Range("B5").Value = 0.9293761783407035
$ws.Range("D5").Value = 0.02949906414341541
$ws.Range("E5").Value = 0.3662447868187053
$ws.Range("F5").Value = 0.8047354716712007
$ws.Range("G5").Value = 0.002443418938322672
$ws.Range("K5").Value = 0.2920327093543449
$ws.Range("L5").Value = 0.09892066805338828
$ws.Range("M5").Value = 0.1828059543812408
$ws.Range("N5").Value = 2.186872375052797
$ws.Range("O5").Value = 2.828009678506589

$ws.Range("B6").Value = 0.9286533945870303
$ws.Range("D6").Value = 0.02942356910596544
$ws.Range("E6").Value = 0.3663546454517999
$ws.Range("F6").Value = 0.8046508078757668
$ws.Range("G6").Value = 0.002443529694098807
$ws.Range("K6").Value = 0.29007218984529
$ws.Range("L6").Value = 0.09868685000395061
$ws.Range("M6").Value = 0.1825929984847896
$ws.Range("N6").Value = 2.187810884230586
$ws.Range("O6").Value = 2.828280489731156

$ws.Range("B7").Value = 0.9337091962230488
$ws.Range("D7").Value = 0.02994658117256677
$ws.Range("E7").Value = 0.3656000219259306
$ws.Range("F7").Value = 0.8052798438729951
$ws.Range("G7").Value = 0.002442768180800969
$ws.Range("K7").Value = 0.3036686288112378
$ws.Range("L7").Value = 0.1003140990544864
$ws.Range("M7").Value = 0.1840790766666487
$ws.Range("N7").Value = 2.181354754490063
$ws.Range("O7").Value = 2.826536125465026

$ws.Range("B8").Value = 0.9569782559154874
$ws.Range("D8").Value = 0.03222949533812169
$ws.Range("E8").Value = 0.3624656134052286
$ws.Range("F8").Value = 0.8090790465833493
$ws.Range("G8").Value = 0.002439586854391587
$ws.Range("K8").Value = 0.363373303193498
$ws.Range("L8").Value = 0.1076005544326222
$ws.Range("M8").Value = 0.1908322995784957
$ws.Range("N8").Value = 2.154303804960953
$ws.Range("O8").Value = 2.822187233672082

$ws.Range("B9").Value = 1.006025252499711
$ws.Range("D9").Value = 0.036615274454185
$ws.Range("E9").Value = 0.3570192063644129
$ws.Range("F9").Value = 0.8202045156827822
$ws.Range("G9").Value = 0.002433987837740301
$ws.Range("K9").Value = 0.479369088098025
$ws.Range("L9").Value = 0.1222649091925518
$ws.Range("M9").Value = 0.2047713249325689
$ws.Range("N9").Value = 2.106423132877587
$ws.Range("O9").Value = 2.825688903273516

$ws.Range("B10").Value = 1.044122063902876
$ws.Range("D10").Value = 0.03978142215520108
$ws.Range("E10").Value = 0.3534419028636648
$ws.Range("F10").Value = 0.8306190833792897
$ws.Range("G10").Value = 0.00243026030799733
$ws.Range("K10").Value = 0.5639129010883153
$ws.Range("L10").Value = 0.1332653279181528
$ws.Range("M10").Value = 0.2154322918691207
$ws.Range("N10").Value = 2.074389218582035
$ws.Range("O10").Value = 2.835596754665545

$ws.Range("B11").Value = 1.061896912700973
$ws.Range("D11").Value = 0.0412094932510243
$ws.Range("E11").Value = 0.3519060328180768
$ws.Range("F11").Value = 0.8358436567501997
$ws.Range("G11").Value = 0.002428647577079894
$ws.Range("K11").Value = 0.6022236506331353
$ws.Range("L11").Value = 0.1383184317620447
$ws.Range("M11").Value = 0.2203725912393466
$ws.Range("N11").Value = 2.060497699795806
$ws.Range("O11").Value = 2.841697517858563

$ws.Range("B12").Value = 1.068691243217501
$ws.Range("D12").Value = 0.04174849229332978
$ws.Range("E12").Value = 0.3513375489736061
$ws.Range("F12").Value = 0.8378920544620598
$ws.Range("G12").Value = 0.002428048743196537
$ws.Range("K12").Value = 0.6167091070843185
$ws.Range("L12").Value = 0.1402388848474772
$ws.Range("M12").Value = 0.2222562795641778
$ws.Range("N12").Value = 2.055335204824694
$ws.Range("O12").Value = 2.844236839944926

$ws.Range("B13").Value = 1.067225152062264
$ws.Range("D13").Value = 0.04163248875253345
$ws.Range("E13").Value = 0.3514593993302393
$ws.Range("F13").Value = 0.8374477845383836
$ws.Range("G13").Value = 0.002428177185515653
$ws.Range("K13").Value = 0.6135903928279447
$ws.Range("L13").Value = 0.1398249731085457
$ws.Range("M13").Value = 0.2218500215758255
$ws.Range("N13").Value = 2.056442685146614
$ws.Range("O13").Value = 2.843679763498329

$ws.Range("B14").Value = 1.062454619124225
$ws.Range("D14").Value = 0.04125387283956172
$ws.Range("E14").Value = 0.3518590006321123
$ws.Range("F14").Value = 0.836010777777858
$ws.Range("G14").Value = 0.002428598073014777
$ws.Range("K14").Value = 0.603415822159775
$ws.Range("L14").Value = 0.1384762897617406
$ws.Range("M14").Value = 0.2205273054805659
$ws.Range("N14").Value = 2.06007101564328
$ws.Range("O14").Value = 2.841901838389475

$ws.Range("B15").Value = 1.05954076842545
$ws.Range("D15").Value = 0.04102172711957763
$ws.Range("E15").Value = 0.3521054751143602
$ws.Range("F15").Value = 0.8351396795113857
$ws.Range("G15").Value = 0.002428857423141862
$ws.Range("K15").Value = 0.5971807213963984
$ws.Range("L15").Value = 0.137651085057783
$ws.Range("M15").Value = 0.2197187806198286
$ws.Range("N15").Value = 2.062306225583788
$ws.Range("O15").Value = 2.840842640904356

$ws.Range("B16").Value = 1.042969333628093
$ws.Range("D16").Value = 0.03968784642346179
$ws.Range("E16").Value = 0.3535441130178762
$ws.Range("F16").Value = 0.8302874388747341
$ws.Range("G16").Value = 0.002430367368155407
$ws.Range("K16").Value = 0.5614061595036901
$ws.Range("L16").Value = 0.1329360733622025
$ws.Range("M16").Value = 0.2151112439758904
$ws.Range("N16").Value = 2.075310763929693
$ws.Range("O16").Value = 2.835230117712172

$ws.Range("B17").Value = 1.032916747336202
$ws.Range("D17").Value = 0.038866406630369
$ws.Range("E17").Value = 0.3544500713730603
$ws.Range("F17").Value = 0.8274354228910283
$ws.Range("G17").Value = 0.0024313148748796
$ws.Range("K17").Value = 0.5394211119965178
$ws.Range("L17").Value = 0.1300560433830924
$ws.Range("M17").Value = 0.2123077818458299
$ws.Range("N17").Value = 2.083463022716963
$ws.Range("O17").Value = 2.832195136820616

$ws.Range("B18").Value = 1.027176627783092
$ws.Range("D18").Value = 0.03839278682696801
$ws.Range("E18").Value = 0.3549797657202696
$ws.Range("F18").Value = 0.8258408546671632
$ws.Range("G18").Value = 0.002431867665574925
$ws.Range("K18").Value = 0.5267619389464926
$ws.Range("L18").Value = 0.1284041397558155
$ws.Range("M18").Value = 0.2107038368302128
$ws.Range("N18").Value = 2.088216072019328
$ws.Range("O18").Value = 2.830599500344704

$ws.Range("B19").Value = 1.025240327818125
$ws.Range("D19").Value = 0.03823223058038394
$ws.Range("E19").Value = 0.3551605911683442
$ws.Range("F19").Value = 0.8253088353029199
$ws.Range("G19").Value = 0.002432056174164483
$ws.Range("K19").Value = 0.5224733824183261
$ws.Range("L19").Value = 0.1278456286053142
$ws.Range("M19").Value = 0.2101622376829582
$ws.Range("N19").Value = 2.089836379257008
$ws.Range("O19").Value = 2.830085010644638

$ws.Range("B20").Value = 1.03398253353518
$ws.Range("D20").Value = 0.03895396941596374
$ws.Range("E20").Value = 0.354352739605214
$ws.Range("F20").Value = 0.8277342814111108
$ws.Range("G20").Value = 0.002431213203263065
$ws.Range("K20").Value = 0.5417629066102734
$ws.Range("L20").Value = 0.1303621506525587
$ws.Range("M20").Value = 0.2126053329880406
$ws.Range("N20").Value = 2.082588568135069
$ws.Range("O20").Value = 2.83250269131571

$ws.Range("B21").Value = 1.063854124336558
$ws.Range("D21").Value = 0.0413651300130482
$ws.Range("E21").Value = 0.35174127233861
$ws.Range("F21").Value = 0.8364309631093789
$ws.Range("G21").Value = 0.002428474126305847
$ws.Range("K21").Value = 0.6064049420003812
$ws.Range("L21").Value = 0.1388722427704749
$ws.Range("M21").Value = 0.2209154700024385
$ws.Range("N21").Value = 2.059002629336018
$ws.Range("O21").Value = 2.842417841046938

$ws.Range("B22").Value = 1.083746124666931
$ws.Range("D22").Value = 0.04293057430216862
$ws.Range("E22").Value = 0.350110966020905
$ws.Range("F22").Value = 0.8425225313804816
$ws.Range("G22").Value = 0.002426753157731986
$ws.Range("K22").Value = 0.6485237704813187
$ws.Range("L22").Value = 0.1444745624683179
$ws.Range("M22").Value = 0.2264217845145353
$ws.Range("N22").Value = 2.044158625773552
$ws.Range("O22").Value = 2.850233276122452

$ws.Range("B23").Value = 1.073095782735777
$ws.Range("D23").Value = 0.04209602523002332
$ws.Range("E23").Value = 0.3509741087197114
$ws.Range("F23").Value = 0.8392340530689069
$ws.Range("G23").Value = 0.002427665359266582
$ws.Range("K23").Value = 0.6260561187700944
$ws.Range("L23").Value = 0.1414808239320564
$ws.Range("M23").Value = 0.2234761234129223
$ws.Range("N23").Value = 2.052028919079322
$ws.Range("O23").Value = 2.845939872751273

$ws.Range("B24").Value = 1.033500569307932
$ws.Range("D24").Value = 0.03891438652473056
$ws.Range("E24").Value = 0.3543967157386172
$ws.Range("F24").Value = 0.8275990270059737
$ws.Range("G24").Value = 0.002431259143803831
$ws.Range("K24").Value = 0.5407042425100599
$ws.Range("L24").Value = 0.1302237475143784
$ws.Range("M24").Value = 0.2124707858033048
$ws.Range("N24").Value = 2.082983702819131
$ws.Range("O24").Value = 2.83236318112759

$ws.Range("B25").Value = 0.9923931333228211
$ws.Range("D25").Value = 0.03543861014080107
$ws.Range("E25").Value = 0.3584179232799292
$ws.Range("F25").Value = 0.8168013018343743
$ws.Range("G25").Value = 0.002435434448161532
$ws.Range("K25").Value = 0.448107045908074
$ws.Range("L25").Value = 0.1182578537988093
$ws.Range("M25").Value = 0.2009263864045252
$ws.Range("N25").Value = 2.118823621447213
$ws.Range("O25").Value = 2.823453577581205
